# natmiOut/OldD7/LR-pairs_lrc2p/Cd80-Cd274.xlsx -- "Natmi following Dr Hou advice"
#
# The sending/target cluster set gained a 4th cluster ("ECs"), turning the
# original 3x4=12 pair table (rows 2-13) into a 4x4=16 pair table
# (rows 2-17), and all of the underlying ligand/receptor expression
# statistics were recomputed accordingly. Columns A-D are
# Sending cluster / Ligand symbol / Receptor symbol / Target cluster;
# columns E-T are the recomputed numeric statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Cd80", "Cd274", "ECs", 2, 0.6666666666666666, 0.539733, 1.619199, 0.0338608224813917, 0.0338608224813917, 3, 1, 12.230514, 36.691542, 0.4767524845277322, 0.4767524845277322, 6.601212012762, 59.410908114858, 0.01614323124615598, 0.01614323124615598),
    @("ECs", "Cd80", "Cd274", "FAPs", 2, 0.6666666666666666, 0.539733, 1.619199, 0.0338608224813917, 0.0338608224813917, 3, 1, 2.417539666666667, 7.252619, 0.09423708950643275, 0.09423708950643275, 1.304825936909, 11.743433432181, 0.00319094535894034, 0.00319094535894034),
    @("ECs", "Cd80", "Cd274", "M2", 2, 0.6666666666666666, 0.539733, 1.619199, 0.0338608224813917, 0.0338608224813917, 3, 1, 10.06935333333333, 30.20806, 0.3925091962001163, 0.3925091962001163, 5.43476228266, 48.91286054394, 0.01329068421484588, 0.01329068421484588),
    @("ECs", "Cd80", "Cd274", "sCs", 2, 0.6666666666666666, 0.539733, 1.619199, 0.0338608224813917, 0.0338608224813917, 3, 1, 0.9363953333333334, 2.809186, 0.03650122976571881, 0.03650122976571881, 0.505403462446, 4.548631162014, 0.001235961661449495, 0.001235961661449495),
    @("FAPs", "Cd80", "Cd274", "ECs", 3, 1, 3.713472666666666, 11.140418, 0.2329693362375475, 0.2329693362375475, 3, 1, 12.230514, 36.691542, 0.4767524845277322, 0.4767524845277322, 45.41767943828399, 408.7591149445559, 0.1110687098700274, 0.1110687098700274),
    @("FAPs", "Cd80", "Cd274", "FAPs", 3, 1, 3.713472666666666, 11.140418, 0.2329693362375475, 0.2329693362375475, 3, 1, 2.417539666666667, 7.252619, 0.09423708950643275, 0.09423708950643275, 8.97746747274911, 80.79720725474199, 0.02195435219127199, 0.02195435219127199),
    @("FAPs", "Cd80", "Cd274", "M2", 3, 1, 3.713472666666666, 11.140418, 0.2329693362375475, 0.2329693362375475, 3, 1, 10.06935333333333, 30.20806, 0.3925091962001163, 0.3925091962001163, 37.39226837434222, 336.5304153690799, 0.09144260690587439, 0.09144260690587439),
    @("FAPs", "Cd80", "Cd274", "sCs", 3, 1, 3.713472666666666, 11.140418, 0.2329693362375475, 0.2329693362375475, 3, 1, 0.9363953333333334, 2.809186, 0.03650122976571881, 0.03650122976571881, 3.477278475527555, 31.295506279748, 0.008503667270373723, 0.008503667270373723),
    @("M2", "Cd80", "Cd274", "ECs", 3, 1, 10.63672166666667, 31.910165, 0.6673079914308979, 0.6673079914308979, 3, 1, 12.230514, 36.691542, 0.4767524845277322, 0.4767524845277322, 130.09257325827, 1170.83315932443, 0.3181407428598912, 0.3181407428598912),
    @("M2", "Cd80", "Cd274", "FAPs", 3, 1, 10.63672166666667, 31.910165, 0.6673079914308979, 0.6673079914308979, 3, 1, 2.417539666666667, 7.252619, 0.09423708950643275, 0.09423708950643275, 25.71469655245945, 231.432268972135, 0.06288516291683138, 0.06288516291683138),
    @("M2", "Cd80", "Cd274", "M2", 3, 1, 10.63672166666667, 31.910165, 0.6673079914308979, 0.6673079914308979, 3, 1, 10.06935333333333, 30.20806, 0.3925091962001163, 0.3925091962001163, 107.1049087699889, 963.9441789299, 0.2619245233344558, 0.2619245233344558),
    @("M2", "Cd80", "Cd274", "sCs", 3, 1, 10.63672166666667, 31.910165, 0.6673079914308979, 0.6673079914308979, 3, 1, 0.9363953333333334, 2.809186, 0.03650122976571881, 0.03650122976571881, 9.960176530632221, 89.64158877569, 0.02435756231971952, 0.02435756231971952),
    @("sCs", "Cd80", "Cd274", "ECs", 3, 1, 1.049821333333333, 3.149464, 0.06586184985016284, 0.06586184985016284, 3, 1, 12.230514, 36.691542, 0.4767524845277322, 0.4767524845277322, 12.839854514832, 115.558690633488, 0.03139980055165758, 0.03139980055165759),
    @("sCs", "Cd80", "Cd274", "FAPs", 3, 1, 1.049821333333333, 3.149464, 0.06586184985016284, 0.06586184985016284, 3, 1, 2.417539666666667, 7.252619, 0.09423708950643275, 0.09423708950643275, 2.537984716246222, 22.841862446216, 0.00620662903938903, 0.00620662903938903),
    @("sCs", "Cd80", "Cd274", "M2", 3, 1, 1.049821333333333, 3.149464, 0.06586184985016284, 0.06586184985016284, 3, 1, 10.06935333333333, 30.20806, 0.3925091962001163, 0.3925091962001163, 10.57102194220444, 95.13919747984, 0.02585138174494016, 0.02585138174494016),
    @("sCs", "Cd80", "Cd274", "sCs", 3, 1, 1.049821333333333, 3.149464, 0.06586184985016284, 0.06586184985016284, 3, 1, 0.9363953333333334, 2.809186, 0.03650122976571881, 0.03650122976571881, 0.9830477973671111, 8.847430176304, 0.002404038514176067, 0.002404038514176067)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    $ws.Cells.Item($r, 1).Value2  = $values[0]   # A: Sending cluster
    $ws.Cells.Item($r, 2).Value2  = $values[1]   # B: Ligand symbol
    $ws.Cells.Item($r, 3).Value2  = $values[2]   # C: Receptor symbol
    $ws.Cells.Item($r, 4).Value2  = $values[3]   # D: Target cluster
    $ws.Cells.Item($r, 5).Value2  = $values[4]   # E
    $ws.Cells.Item($r, 6).Value2  = $values[5]   # F
    $ws.Cells.Item($r, 7).Value2  = $values[6]   # G
    $ws.Cells.Item($r, 8).Value2  = $values[7]   # H
    $ws.Cells.Item($r, 9).Value2  = $values[8]   # I
    $ws.Cells.Item($r, 10).Value2 = $values[9]   # J
    $ws.Cells.Item($r, 11).Value2 = $values[10]  # K
    $ws.Cells.Item($r, 12).Value2 = $values[11]  # L
    $ws.Cells.Item($r, 13).Value2 = $values[12]  # M
    $ws.Cells.Item($r, 14).Value2 = $values[13]  # N
    $ws.Cells.Item($r, 15).Value2 = $values[14]  # O
    $ws.Cells.Item($r, 16).Value2 = $values[15]  # P
    $ws.Cells.Item($r, 17).Value2 = $values[16]  # Q
    $ws.Cells.Item($r, 18).Value2 = $values[17]  # R
    $ws.Cells.Item($r, 19).Value2 = $values[18]  # S
    $ws.Cells.Item($r, 20).Value2 = $values[19]  # T
}
